# Generate Report for Handback
#
# The existing handed-back file 7e2f4365-c3c1-442c-a624-bbffcb321438.md was
# regenerated as 6c864454-4b85-46d2-be32-8d1575f62b92.md (row 2 on every
# sheet gets its identifiers / timestamps refreshed), and a brand new file
# 9c93726a-a17b-4071-91f6-485915fec7c8.md was handed back too (becomes the
# new row 3 on every sheet).
#
# NOTE: a leading "'" forces literal text so values like "True"/"False"/""
# are not auto-coerced into booleans / dropped, and Range.Hyperlinks.Delete()
# clears every hyperlink on the sheet (not just the target range) so all
# hyperlinks on a sheet are deleted once up front and re-added afterwards.

$wb = $excel.ActiveWorkbook

$oldGuid = "7e2f4365-c3c1-442c-a624-bbffcb321438"
$guidA = "6c864454-4b85-46d2-be32-8d1575f62b92"
$guidB = "9c93726a-a17b-4071-91f6-485915fec7c8"

$hashA = "be0966b3a36dfc5d82d17d618880c3c24567f6f2"
$hashB = "c71666dde1fbf2c7e6c1ab9f67e2953b9dc571e5"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item("Overview")
$lo.ListRows.Add() | Out-Null

# Wipe the sheet's hyperlinks once (Range.Hyperlinks.Delete is sheet-wide
# in this engine) then rebuild them all after the cell values are set.
$ws.Range("A1").Hyperlinks.Delete()

# Row 2 - refresh the existing file's identifiers
$ws.Range("A2").Value = "'$guidA.md"
$ws.Range("C2").Value = "'.md"
$ws.Range("E2").Value = "'Handed back: in sync with en-US"
$ws.Range("F2").Value = "'Handed back: in sync with en-US"
$ws.Range("G2").Value = "'2016-08-17 21:02:15"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 - the newly handed back file
$ws.Range("A3").Value = "'$guidB.md"
$ws.Range("C3").Value = "'.md"
$ws.Range("E3").Value = "'Handed back: in sync with en-US"
$ws.Range("F3").Value = "'Handed back: in sync with en-US"
$ws.Range("G3").Value = "'2016-08-17 21:02:15"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidA.md", "", "", "e2e\$guidA.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidB.md", "", "", "e2e\$guidB.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item("zh_cn")
$lo.ListRows.Add() | Out-Null

$ws.Range("A1").Hyperlinks.Delete()

# Row 2 - refresh the existing file's identifiers
$ws.Range("B2").Value = "'.md"
$ws.Range("C2").Value = "'Handed back: in sync with en-US"
$ws.Range("D2").Value = "'e2e"
$ws.Range("E2").Value = "'ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "'$guidA.$hashA.zh-cn.xlf"
$ws.Range("H2").Value = "'2016-08-17 21:02:05"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J2").Value = "'$guidA.$hashA.zh-cn.xlf"
$ws.Range("K2").Value = "'2016-08-17 21:02:31"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = "'"

# Row 3 - the newly handed back file
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Handed back: in sync with en-US"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "'$guidB.$hashB.zh-cn.xlf"
$ws.Range("H3").Value = "'2016-08-17 21:02:05"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "'$guidB.$hashB.zh-cn.xlf"
$ws.Range("K3").Value = "'2016-08-17 21:02:31"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidA.md", "", "", "$guidA.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/70bc352fab81dd2156f452ef32d2fd13f6631e5c/e2e/$guidA.md", "", "", "$guidA.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidB.md", "", "", "$guidB.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/70bc352fab81dd2156f452ef32d2fd13f6631e5c/e2e/$guidB.md", "", "", "$guidB.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item("de_de")
$lo.ListRows.Add() | Out-Null

$ws.Range("A1").Hyperlinks.Delete()

# Row 2 - refresh the existing file's identifiers
$ws.Range("B2").Value = "'.md"
$ws.Range("C2").Value = "'Handed back: in sync with en-US"
$ws.Range("D2").Value = "'e2e"
$ws.Range("E2").Value = "'ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "'$guidA.$hashA.de-de.xlf"
$ws.Range("H2").Value = "'2016-08-17 21:02:15"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J2").Value = "'$guidA.$hashA.de-de.xlf"
$ws.Range("K2").Value = "'2016-08-17 21:02:38"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = "'"

# Row 3 - the newly handed back file
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Handed back: in sync with en-US"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "'$guidB.$hashB.de-de.xlf"
$ws.Range("H3").Value = "'2016-08-17 21:02:15"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "'$guidB.$hashB.de-de.xlf"
$ws.Range("K3").Value = "'2016-08-17 21:02:38"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidA.md", "", "", "$guidA.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f7d5b13548a9c284c19fda70560cf0000846d9e/e2e/$guidA.md", "", "", "$guidA.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0498d88d3b6c04ea7567169e2b221687fcb03ec2/e2e/$guidB.md", "", "", "$guidB.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f7d5b13548a9c284c19fda70560cf0000846d9e/e2e/$guidB.md", "", "", "$guidB.md") | Out-Null
